# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" -> row number : new F value
$sheetExhibition = $wb.Worksheets.Item("展览")
$exhibitionUpdates = @{
    3  = 495
    4  = 1278
    5  = 1133
    6  = 14208
    7  = 15997
    9  = 70
    22 = 68
    24 = 6366
    26 = 1108
    27 = 5638
    28 = 85
    31 = 4644
}
foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型" -> row number : new F value
$sheetAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    3  = 495
    4  = 1278
    5  = 1133
    6  = 14208
    7  = 15997
    9  = 70
    22 = 68
    25 = 6366
    27 = 1108
    29 = 5638
    30 = 85
    33 = 4644
}
foreach ($row in $allUpdates.Keys) {
    $sheetAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}

$wb.Save()
